{"js": "const replacements = [\n  [\"2024-02-08 Thursday\", \"2024-02-09 Friday\"],\n  [\"326\u00f78=\", \"873\u00f73=\"],\n  [\"617\u00f73=\", \"721\u00f79=\"],\n  [\"333\u00f79=\", \"255\u00f74=\"],\n  [\"117\u00f77=\", \"589\u00f75=\"],\n  [\"400\u00f78=\", \"586\u00f78=\"],\n  [\"120\u00f73=\", \"244\u00f77=\"],\n  [\"987\u00f79=\", \"541\u00f74=\"],\n  [\"491\u00f79=\", \"777\u00f78=\"],\n  [\"114\u00f74=\", \"434\u00f77=\"],\n  [\"764\u00f78=\", \"877\u00f77=\"],\n  [\"337\u00f75=\", \"881\u00f78=\"],\n  [\"592\u00f78=\", \"464\u00f74=\"],\n  [\"123\u00f72=\", \"433\u00f73=\"],\n  [\"710\u00f75=\", \"798\u00f74=\"],\n  [\"681\u00f75=\", \"803\u00f73=\"],\n  [\"893\u00f73=\", \"675\u00f72=\"],\n  [\"413\u00f73=\", \"139\u00f72=\"],\n  [\"661\u00f74=\", \"284\u00f76=\"],\n  [\"989\u00f75=\", \"230\u00f73=\"],\n  [\"282\u00f72=\", \"739\u00f79=\"],\n  [\"808\u00f76=\", \"381\u00f76=\"],\n  [\"991\u00f76=\", \"369\u00f75=\"],\n  [\"775\u00f75=\", \"401\u00f78=\"],\n  [\"352\u00f76=\", \"981\u00f78=\"],\n  [\"310\u00f75=\", \"986\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-08 Thursday\", \"2024-02-09 Friday\"),\n    @(\"326\u00f78=\", \"873\u00f73=\"),\n    @(\"617\u00f73=\", \"721\u00f79=\"),\n    @(\"333\u00f79=\", \"255\u00f74=\"),\n    @(\"117\u00f77=\", \"589\u00f75=\"),\n    @(\"400\u00f78=\", \"586\u00f78=\"),\n    @(\"120\u00f73=\", \"244\u00f77=\"),\n    @(\"987\u00f79=\", \"541\u00f74=\"),\n    @(\"491\u00f79=\", \"777\u00f78=\"),\n    @(\"114\u00f74=\", \"434\u00f77=\"),\n    @(\"764\u00f78=\", \"877\u00f77=\"),\n    @(\"337\u00f75=\", \"881\u00f78=\"),\n    @(\"592\u00f78=\", \"464\u00f74=\"),\n    @(\"123\u00f72=\", \"433\u00f73=\"),\n    @(\"710\u00f75=\", \"798\u00f74=\"),\n    @(\"681\u00f75=\", \"803\u00f73=\"),\n    @(\"893\u00f73=\", \"675\u00f72=\"),\n    @(\"413\u00f73=\", \"139\u00f72=\"),\n    @(\"661\u00f74=\", \"284\u00f76=\"),\n    @(\"989\u00f75=\", \"230\u00f73=\"),\n    @(\"282\u00f72=\", \"739\u00f79=\"),\n    @(\"808\u00f76=\", \"381\u00f76=\"),\n    @(\"991\u00f76=\", \"369\u00f75=\"),\n    @(\"775\u00f75=\", \"401\u00f78=\"),\n    @(\"352\u00f76=\", \"981\u00f78=\"),\n    @(\"310\u00f75=\", \"986\u00f73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}"}
